$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for new column F, matching style of existing header cells (E1 etc.)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Timestamp values for each data row
$timestamps = @(
    "2021-10-05 10:50:30.808023",
    "2021-10-05 10:50:30.808035",
    "2021-10-05 10:50:30.808039",
    "2021-10-05 10:50:30.808042",
    "2021-10-05 10:50:30.808046",
    "2021-10-05 10:50:30.808049",
    "2021-10-05 10:50:30.808052",
    "2021-10-05 10:50:30.808055",
    "2021-10-05 10:50:30.808058",
    "2021-10-05 10:50:30.808061",
    "2021-10-05 10:50:30.808064",
    "2021-10-05 10:50:30.808067",
    "2021-10-05 10:50:30.808070"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
